$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 830
$ws.Range("I2").Value = 700.1429000000001
$ws.Range("K2").Value = 700.1429000000001
$ws.Range("M2").Value = -587.1429000000001
$ws.Range("H4").Value = 2740
$ws.Range("I4").Value = 2740
$ws.Range("K4").Value = 2740
$ws.Range("M4").Value = -2626
$ws.Range("H5").Value = 106.4
$ws.Range("I5").Value = 120.5
$ws.Range("K5").Value = 120.5
$ws.Range("M5").Value = -5.5
$ws.Range("H9").Value = 338.69232
$ws.Range("I9").Value = 429.4
$ws.Range("K9").Value = 429.4
$ws.Range("M9").Value = -260.4
$ws.Range("H62").Value = 3002
$ws.Range("I62").Value = 3002
$ws.Range("K62").Value = 3002
$ws.Range("M62").Value = -2378
$ws.Range("H65").Value = 3002
$ws.Range("I65").Value = 3002
$ws.Range("K65").Value = 15010
$ws.Range("M65").Value = -11890
$ws.Range("H70").Value = 1133.5555
$ws.Range("I70").Value = 850
$ws.Range("J70").Value = 1214.5714
$ws.Range("K70").Value = 2550
$ws.Range("L70").Value = 3643.7142
$ws.Range("M70").Value = -2280
$ws.Range("N70").Value = -4183.7142
$ws.Range("H73").Value = 1133.5555
$ws.Range("I73").Value = 850
$ws.Range("J73").Value = 1214.5714
$ws.Range("K73").Value = 2550
$ws.Range("L73").Value = 3643.7142
$ws.Range("M73").Value = -1614
$ws.Range("N73").Value = -5515.7142
$ws.Range("H74").Value = 2833.3333
$ws.Range("I74").Value = 2833.3333
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2833.3333
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1897.3333
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2833.3333
$ws.Range("I77").Value = 2833.3333
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 14166.6665
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -9486.666499999999
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 2041.8889
$ws.Range("I132").Value = 948.1667
$ws.Range("J132").Value = 4229.3335
$ws.Range("K132").Value = 2844.5001
$ws.Range("L132").Value = 12688.0005
$ws.Range("M132").Value = -314.5001000000002
$ws.Range("N132").Value = -17748.0005
$ws.Range("H135").Value = 1665.7858
$ws.Range("I135").Value = 1577.75
$ws.Range("K135").Value = 14199.75
$ws.Range("M135").Value = -11664.75
$ws.Range("H137").Value = 2112
$ws.Range("I137").Value = 2066
$ws.Range("K137").Value = 6198
$ws.Range("M137").Value = -3648
$ws.Range("H138").Value = 3773.7144
$ws.Range("I138").Value = 845
$ws.Range("K138").Value = 2535
$ws.Range("M138").Value = 2605

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3222.4285
$ws.Range("J61").Value = 6999
$ws.Range("L61").Value = 6999
$ws.Range("N61").Value = -7423
$ws.Range("H74").Value = 1719.4286
$ws.Range("I74").Value = 1357.3
$ws.Range("K74").Value = 1357.3
$ws.Range("M74").Value = -483.3
$ws.Range("H77").Value = 1719.4286
$ws.Range("I77").Value = 1357.3
$ws.Range("K77").Value = 6786.5
$ws.Range("M77").Value = -2418.5
$ws.Range("H132").Value = 2146.647
$ws.Range("I132").Value = 2032.9333
$ws.Range("K132").Value = 6098.7999
$ws.Range("M132").Value = -3568.7999
$ws.Range("H136").Value = 3222.4285
$ws.Range("J136").Value = 6999
$ws.Range("L136").Value = 20997
$ws.Range("N136").Value = -26097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2981.077
$ws.Range("I105").Value = 2768.5454
$ws.Range("J105").Value = 4150
$ws.Range("K105").Value = 2768.5454
$ws.Range("L105").Value = 4150
$ws.Range("M105").Value = -1021.5454
$ws.Range("N105").Value = -7644
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 7765.7856
$ws.Range("I134").Value = 3833.75
$ws.Range("J134").Value = 13008.5
$ws.Range("K134").Value = 11501.25
$ws.Range("L134").Value = 39025.5
$ws.Range("M134").Value = -8966.25
$ws.Range("N134").Value = -44095.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 26194.445
$ws.Range("J50").Value = 26194.445
$ws.Range("L50").Value = 26194.445
$ws.Range("N50").Value = -27444.445
$ws.Range("H51").Value = 21166.666
$ws.Range("J51").Value = 21400
$ws.Range("L51").Value = 21400
$ws.Range("N51").Value = -22872
$ws.Range("H58").Value = 6870
$ws.Range("I58").Value = 6870
$ws.Range("K58").Value = 6870
$ws.Range("M58").Value = -6667
$ws.Range("H59").Value = 30000
$ws.Range("H60").Value = 24166.666
$ws.Range("J60").Value = 28333.334
$ws.Range("L60").Value = 28333.334
$ws.Range("N60").Value = -29355.334
$ws.Range("H61").Value = 21166.666
$ws.Range("J61").Value = 21400
$ws.Range("L61").Value = 21400
$ws.Range("N61").Value = -22096
$ws.Range("H96").Value = 9184.799999999999
$ws.Range("J96").Value = 9184.799999999999
$ws.Range("L96").Value = 9184.799999999999
$ws.Range("N96").Value = -14676.8
$ws.Range("H125").Value = 45000
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -49920
$ws.Range("H132").Value = 976.375
$ws.Range("I132").Value = 830.1429000000001
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2490.4287
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = 39.57129999999961
$ws.Range("N132").Value = -11060
$ws.Range("H134").Value = 3259.3333
$ws.Range("I134").Value = 2511.2
$ws.Range("K134").Value = 7533.599999999999
$ws.Range("M134").Value = -4998.599999999999
$ws.Range("H136").Value = 6870
$ws.Range("I136").Value = 6870
$ws.Range("K136").Value = 20610
$ws.Range("M136").Value = -18060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H4").Value = 2999999.5
$ws.Range("I4").Value = 5000000
$ws.Range("K4").Value = 15000000
$ws.Range("M4").Value = -14999888
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H46").Value = 949
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 949
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2847
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3029
$ws.Range("H103").Value = 2382.3
$ws.Range("J103").Value = 2170.5
$ws.Range("L103").Value = 6511.5
$ws.Range("N103").Value = -8269.5
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5000
$ws.Range("I12").Value = 5000
$ws.Range("K12").Value = 5000
$ws.Range("M12").Value = -4860
$ws.Range("H80").Value = 5025.1816
$ws.Range("I80").Value = 5377.7
$ws.Range("J80").Value = 1500
$ws.Range("K80").Value = 5377.7
$ws.Range("L80").Value = 1500
$ws.Range("M80").Value = -4379.7
$ws.Range("N80").Value = -3496
$ws.Range("H83").Value = 5025.1816
$ws.Range("I83").Value = 5377.7
$ws.Range("J83").Value = 1500
$ws.Range("K83").Value = 26888.5
$ws.Range("L83").Value = 7500
$ws.Range("M83").Value = -21896.5
$ws.Range("N83").Value = -17484
$ws.Range("H122").Value = 14666.667
$ws.Range("I122").Value = 14666.667
$ws.Range("K122").Value = 44000.001
$ws.Range("M122").Value = -41550.001
$ws.Range("H132").Value = 3053.353
$ws.Range("I132").Value = 2382.5386
$ws.Range("J132").Value = 5233.5
$ws.Range("K132").Value = 7147.6158
$ws.Range("L132").Value = 15700.5
$ws.Range("M132").Value = -4617.6158
$ws.Range("N132").Value = -20760.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5650
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H55").Value = 5999.5
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 5999.5
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 5999.5
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -6345.5
$ws.Range("H122").Value = 4475.1665
$ws.Range("I122").Value = 4001.3333
$ws.Range("J122").Value = 4949
$ws.Range("K122").Value = 12003.9999
$ws.Range("L122").Value = 14847
$ws.Range("M122").Value = -9553.999899999999
$ws.Range("N122").Value = -19747
$ws.Range("H136").Value = 3662.3333
$ws.Range("I136").Value = 3662.3333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10986.9999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8436.999899999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 50000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51872
$ws.Range("H77").Value = 50000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159360
$ws.Range("H122").Value = 2357
$ws.Range("I122").Value = 999.8
$ws.Range("J122").Value = 5750
$ws.Range("K122").Value = 2999.4
$ws.Range("L122").Value = 17250
$ws.Range("M122").Value = -549.3999999999996
$ws.Range("N122").Value = -22150
$ws.Range("H132").Value = 3555.6875
$ws.Range("I132").Value = 1635.8182
$ws.Range("J132").Value = 7779.4
$ws.Range("K132").Value = 4907.4546
$ws.Range("L132").Value = 23338.2
$ws.Range("M132").Value = -2377.4546
$ws.Range("N132").Value = -28398.2
